# Update the "Estado de Cuenta" database: the "Periodo Mora" (E16:E25) and
# "Valor Mora" (F16/F25) figures are refreshed with the new account-statement
# data (commit: "Actualiza base de datos EC y agrega parte 1 de nuevos
# estado de cuenta"). The period labels are renumbered in ascending order
# (1901 .. 1910); the outstanding "Valor Mora" amount of 19305 stays tied to
# period 1910, which now sits in row 25, while every other period keeps
# 27578.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1901"
$ws.Range("E17").Value = "1902"
$ws.Range("E18").Value = "1903"
$ws.Range("E19").Value = "1904"
$ws.Range("E20").Value = "1905"
$ws.Range("E21").Value = "1906"
$ws.Range("E22").Value = "1907"
$ws.Range("E23").Value = "1908"
$ws.Range("E24").Value = "1909"
$ws.Range("E25").Value = "1910"

$ws.Range("F16").Value = 27578
$ws.Range("F25").Value = 19305
